# Update Step3_DataPts_* sheets with new First_Noticeable_Increase_Index (C),
# First_Noticeable_Increase_Cumulative_Value (E), and Pulse_Width (G) values.
# These reflect the addition of a configurable zero_before_threshold parameter
# that allows dims before the noise_threshold / First Rise Point to be zeroed.

$wb = $excel.ActiveWorkbook

$sheetData = @{
    "Step3_DataPts_0.5" = @{
        2 = @(88, 0.03245730786473743, 51)
        3 = @(88, 0.02336815654301692, 27)
        4 = @(88, 0.02640150572725969, 40)
        5 = @(87, 0.007875619107061919, 18)
        6 = @(88, 0.01343193811605964, 20)
    }
    "Step3_DataPts_0.7" = @{
        2 = @(88, 0.03245730786473743, 75)
        3 = @(88, 0.02336815654301692, 66)
        4 = @(88, 0.02640150572725969, 73)
        5 = @(87, 0.007875619107061919, 72)
        6 = @(88, 0.01343193811605964, 68)
    }
    "Step3_DataPts_0.8" = @{
        2 = @(88, 0.03245730786473743, 77)
        3 = @(88, 0.02336815654301692, 74)
        4 = @(88, 0.02640150572725969, 77)
        5 = @(87, 0.007875619107061919, 77)
        6 = @(88, 0.01343193811605964, 75)
    }
    "Step3_DataPts_0.9" = @{
        2 = @(88, 0.03245730786473743, 89)
        3 = @(88, 0.02336815654301692, 87)
        4 = @(88, 0.02640150572725969, 88)
        5 = @(87, 0.007875619107061919, 88)
        6 = @(88, 0.01343193811605964, 88)
    }
}

foreach ($sheetName in $sheetData.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $sheetData[$sheetName]
    foreach ($r in $rows.Keys) {
        $vals = $rows[$r]
        $ws.Cells.Item($r, 3).Value = $vals[0]   # Column C
        $ws.Cells.Item($r, 5).Value = $vals[1]   # Column E
        $ws.Cells.Item($r, 7).Value = $vals[2]   # Column G
    }
}
